$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing test string used by B17 ("테스트" -> "테스트1")
$ws.Range("B17").Value = "테스트1"

# Add a new row 18, copying formatting (styles/borders/number formats) from row 17,
# then overwrite with the new row's actual values.
$ws.Range("A17:H17").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)

$ws.Range("A18").Value = "A05"
$ws.Range("B18").Value = "테스트2"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "없음"
$ws.Range("F18").Value = 162000000
$ws.Range("G18").Value = "임종완"
$ws.Range("H18").Value = 45839

# Update the active selection to reflect where the user ended up after the edit.
[void]$ws.Range("G22").Select()
